$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(444).Insert()

$ws.Cells.Item(444, 1).Value = 4
$ws.Cells.Item(444, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(444, 3).Value = "Los Lagos"
$ws.Cells.Item(444, 4).Value = 45127
$ws.Cells.Item(444, 5).Value = 10
$ws.Cells.Item(444, 6).Value = 100112045
$ws.Cells.Item(444, 7).Value = "Zapallo"
$ws.Cells.Item(444, 8).Value = "Paine"
$ws.Cells.Item(444, 9).Value = "1a (guarda)"
$ws.Cells.Item(444, 10).Value = 500
$ws.Cells.Item(444, 11).Value = 550
$ws.Cells.Item(444, 12).Value = 600
$ws.Cells.Item(444, 13).Value = 575
$ws.Cells.Item(444, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(444, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(444, 16).Value = 575
$ws.Cells.Item(444, 17).Value = 1
$ws.Cells.Item(444, 18).Value = "Hortaliza"
